$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.40130033333333
$ws.Range("H2").Value = 109.203901
$ws.Range("I2").Value = 0.1897437225523226
$ws.Range("J2").Value = 0.1897437225523226
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.051093
$ws.Range("N2").Value = 0.153279
$ws.Range("O2").Value = 0.01450579975525089
$ws.Range("P2").Value = 0.01450579975525089
$ws.Range("Q2").Value = 1.859851637931
$ws.Range("R2").Value = 16.738664741379
$ws.Range("S2").Value = 0.002752384444159873
$ws.Range("T2").Value = 0.002752384444159874

# Row 3
$ws.Range("G3").Value = 36.40130033333333
$ws.Range("H3").Value = 109.203901
$ws.Range("I3").Value = 0.1897437225523226
$ws.Range("J3").Value = 0.1897437225523226
$ws.Range("O3").Value = 0.2313022967634575
$ws.Range("P3").Value = 0.2313022967634575
$ws.Range("Q3").Value = 29.65627285300122
$ws.Range("R3").Value = 266.906455677011
$ws.Range("S3").Value = 0.04388815882280045
$ws.Range("T3").Value = 0.04388815882280047

# Row 4
$ws.Range("G4").Value = 36.40130033333333
$ws.Range("H4").Value = 109.203901
$ws.Range("I4").Value = 0.1897437225523226
$ws.Range("J4").Value = 0.1897437225523226
$ws.Range("M4").Value = 2.656449666666667
$ws.Range("N4").Value = 7.969348999999999
$ws.Range("O4").Value = 0.7541919034812916
$ws.Range("P4").Value = 0.7541919034812917
$ws.Range("Q4").Value = 96.69822213671655
$ws.Range("R4").Value = 870.2839992304489
$ws.Range("S4").Value = 0.1431031792853623
$ws.Range("T4").Value = 0.1431031792853623

# Row 5
$ws.Range("I5").Value = 0.6107553255746098
$ws.Range("J5").Value = 0.6107553255746098
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.051093
$ws.Range("N5").Value = 0.153279
$ws.Range("O5").Value = 0.01450579975525089
$ws.Range("P5").Value = 0.01450579975525089
$ws.Range("Q5").Value = 5.986571135874001
$ws.Range("R5").Value = 53.87914022286601
$ws.Range("S5").Value = 0.00885949445223835
$ws.Range("T5").Value = 0.008859494452238352

# Row 6
$ws.Range("I6").Value = 0.6107553255746098
$ws.Range("J6").Value = 0.6107553255746098
$ws.Range("O6").Value = 0.2313022967634575
$ws.Range("P6").Value = 0.2313022967634575
$ws.Range("Q6").Value = 95.45889760157712
$ws.Range("R6").Value = 859.130078414194
$ws.Range("S6").Value = 0.1412691095659205
$ws.Range("T6").Value = 0.1412691095659205

# Row 7
$ws.Range("I7").Value = 0.6107553255746098
$ws.Range("J7").Value = 0.6107553255746098
$ws.Range("M7").Value = 2.656449666666667
$ws.Range("N7").Value = 7.969348999999999
$ws.Range("O7").Value = 0.7541919034812916
$ws.Range("P7").Value = 0.7541919034812917
$ws.Range("Q7").Value = 311.2564323560718
$ws.Range("R7").Value = 2801.307891204646
$ws.Range("S7").Value = 0.4606267215564509
$ws.Range("T7").Value = 0.460626721556451

# Row 8
$ws.Range("G8").Value = 38.27317166666666
$ws.Range("H8").Value = 114.819515
$ws.Range("I8").Value = 0.1995009518730676
$ws.Range("J8").Value = 0.1995009518730676
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.051093
$ws.Range("N8").Value = 0.153279
$ws.Range("O8").Value = 0.01450579975525089
$ws.Range("P8").Value = 0.01450579975525089
$ws.Range("Q8").Value = 1.955491159965
$ws.Range("R8").Value = 17.599420439685
$ws.Range("S8").Value = 0.002893920858852663
$ws.Range("T8").Value = 0.002893920858852664

# Row 9
$ws.Range("G9").Value = 38.27317166666666
$ws.Range("H9").Value = 114.819515
$ws.Range("I9").Value = 0.1995009518730676
$ws.Range("J9").Value = 0.1995009518730676
$ws.Range("O9").Value = 0.2313022967634575
$ws.Range("P9").Value = 0.2313022967634575
$ws.Range("Q9").Value = 31.18129329179611
$ws.Range("R9").Value = 280.631639626165
$ws.Range("S9").Value = 0.04614502837473653
$ws.Range("T9").Value = 0.04614502837473655

# Row 10
$ws.Range("G10").Value = 38.27317166666666
$ws.Range("H10").Value = 114.819515
$ws.Range("I10").Value = 0.1995009518730676
$ws.Range("J10").Value = 0.1995009518730676
$ws.Range("M10").Value = 2.656449666666667
$ws.Range("N10").Value = 7.969348999999999
$ws.Range("O10").Value = 0.7541919034812916
$ws.Range("P10").Value = 0.7541919034812917
$ws.Range("Q10").Value = 101.6707541161928
$ws.Range("R10").Value = 915.0367870457349
$ws.Range("S10").Value = 0.1504620026394784
$ws.Range("T10").Value = 0.1504620026394785
